$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Extra mobile data requests")
$rng = $ws1.Range("C2")
$v = $rng.Validation
$v | Get-Member | Out-String | Write-Host
